$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowData {
    param($row, $sourceRow, $phone, $ddd, $date)
    $ws.Range("A$row").Value = "'" + $phone
    if ($ddd -ne $null) {
        $ws.Range("B$row").Value = "'" + $ddd
    } else {
        $ws.Range("B$row").Value = ""
    }
    $ws.Range("C$row").Value = "'" + $date
    # Re-apply the row's original formatting (font/fill/border/alignment/numberformat)
    # so the freshly-typed text doesn't pick up an incidental "number stored as text"
    # style variant (quote-prefix) - matches the plain style used throughout the column.
    $ws.Range("A$sourceRow`:C$sourceRow").Copy()
    $ws.Range("A$row`:C$row").PasteSpecial(-4122)
    $excel.CutCopyMode = $false
}

# Insert 1 new row before old row 12 (shifts everything below down by 1)
$ws.Rows.Item(12).Insert()
Set-RowData 12 13 "+5511966134418" "11" "2024-09-24"

# Insert 2 new rows before what is now row 15 (old row 14, after the first shift)
$ws.Rows.Item(15).Insert()
$ws.Rows.Item(15).Insert()
Set-RowData 15 17 "+556192771804" "61" "2024-09-18"
Set-RowData 16 17 "+556198454144" "61" "2024-09-18"

# Insert 1 new row before what is now row 18 (old row 15, after prior shifts)
$ws.Rows.Item(18).Insert()
Set-RowData 18 19 "+555199100909" "51" "2024-09-13"

# Insert 1 new row before what is now row 20 (old row 16, after prior shifts)
$ws.Rows.Item(20).Insert()
Set-RowData 20 21 "+5511947261969" "11" "2024-09-10"
